$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix height values that were entered as text (e.g. "6Ft", "5.8ft") to proper numbers,
# and correct a typo (167 -> 5.4).
$ws.Range("A27").Value = 6
$ws.Range("A30").Value = 5.8
$ws.Range("A47").Value = 5.4
